$wb = $excel.ActiveWorkbook

# --- MasterMaps: add a "Type" column (D) describing each map's render style ---
$ws1 = $wb.Worksheets.Item("MasterMaps")

$ws1.Cells.Item(1, 4).Value = "Type"
for ($r = 2; $r -le 22; $r++) {
    if ($r -eq 7) {
        # row 7 is the "Classes" map, which is a discrete (categorical) map
        $ws1.Cells.Item($r, 4).Value = "discrete"
    } else {
        $ws1.Cells.Item($r, 4).Value = "streched"
    }
}

# Move the selection on MasterMaps to D8, matching where the edit left off
[void]$ws1.Range("D8").Select()

# MasterMaps becomes the active sheet/tab (was "O.C.")
$ws1.Activate()

Write-Host "Applied MasterMaps Type column + selection updates"
